$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new credit entry in row 10
$ws.Range("A10").Value = "BronyB34r"
$ws.Range("B10").Value = "discord image"

# Update the selection to match the new state
$ws.Range("A4:B10").Select()
$ws.Application.ActiveCell = $ws.Range("A4")
